$d = $word.ActiveDocument

$replacements = @(
    ,@("2024-10-22 Tuesday", "2024-10-23 Wednesday")
    ,@("20×70=1400", "24×57=1368")
    ,@("75×60=4500", "13×96=1248")
    ,@("54×99=5346", "23×38=874")
    ,@("19×63=1197", "40×59=2360")
    ,@("65×32=2080", "43×96=4128")
    ,@("46×90=4140", "33×14=462")
    ,@("22×33=726", "87×96=8352")
    ,@("56×56=3136", "94×56=5264")
    ,@("27×24=648", "89×89=7921")
    ,@("61×48=2928", "97×94=9118")
    ,@("44×63=2772", "26×53=1378")
    ,@("65×73=4745", "47×31=1457")
    ,@("55×25=1375", "55×31=1705")
    ,@("84×19=1596", "19×72=1368")
    ,@("73×61=4453", "15×77=1155")
    ,@("51×63=3213", "71×29=2059")
    ,@("71×12=852", "72×32=2304")
    ,@("38×48=1824", "62×29=1798")
    ,@("49×95=4655", "14×90=1260")
    ,@("20×39=780", "48×11=528")
    ,@("32×63=2016", "74×83=6142")
    ,@("57×92=5244", "67×91=6097")
    ,@("47×26=1222", "50×49=2450")
    ,@("31×72=2232", "44×30=1320")
    ,@("29×41=1189", "95×16=1520")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done."
